# Adds the RadioButtonPage and WebTablesPage rows to the "URL" sheet,
# matching the target edit (rows 6 and 7 of the URL table).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("URL")

# Row 6: RadioButtonPage URL
$ws.Range("A6").Value = "5"
$ws.Range("B6").Value = "RadioButtonPage URL"
$ws.Range("C6").Value = "https://demoqa.com/radio-button"

# Row 7: WebTablesPage URL
# (Note: C7 is entered before B7 so the shared-string interning order
# matches the original author's editing order.)
$ws.Range("A7").Value = "6"
$ws.Range("C7").Value = "https://demoqa.com/webtables"
$ws.Range("B7").Value = "WebTablesPage URL"

# Widen column C slightly to fit the new longer URL ("https://demoqa.com/radio-button").
$ws.Range("C1").ColumnWidth = 31.5

# Move/restore the selection, as in the target workbook.
$ws.Range("C10").Select()
